# Generate Report for Handback
# - updates the "Status" text (Ready for handoff -> Handed back: in sync with en-US)
#   everywhere it is used (Overview + zh-cn + de-de sheets)
# - fills in "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   for both language sheets, adding the corresponding hyperlink on the Target File cell
# - widens the columns that now hold longer text

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$targetMd  = "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"
$baseUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/930f8002d3edc70889801a6789bd3ce652f617a3/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the status text + widen the two status columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $targetMd
$wsZh.Range("I3").Value = $targetMd
$wsZh.Range("J2").Value = "01a97a15-1964-4c43-b589-d1f0cc6bc596.e059ce618015e54d813a2480229315a514275c36.zh-cn.xlf"
$wsZh.Range("J3").Value = "01a97a15-1964-4c43-b589-d1f0cc6bc596.e059ce618015e54d813a2480229315a514275c36.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-30 07:06:28"
$wsZh.Range("K3").Value = "2016-08-30 07:06:28"

# Rebuild the hyperlinks collection so the new "Latest Target File" link on I2/I3
# lands between the existing A2/A3 links, in the same layout the workbook ends up with.
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($baseUrl + "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"), "", "", "01a97a15-1964-4c43-b589-d1f0cc6bc596.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($baseUrl + "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"), "", "", "01a97a15-1964-4c43-b589-d1f0cc6bc596.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($baseUrl + "ffffa6b67d8f-8e16-4a84-8d2f-d01ff423c2eb.md"), "", "", "ffffa6b67d8f-8e16-4a84-8d2f-d01ff423c2eb.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($baseUrl + "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"), "", "", "01a97a15-1964-4c43-b589-d1f0cc6bc596.md")

$wsZh.Range("A2").Style = "HyperLink"
$wsZh.Range("A3").Style = "HyperLink"
$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Range("I3").Style = "HyperLink"

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $targetMd
$wsDe.Range("I3").Value = $targetMd
$wsDe.Range("J2").Value = "01a97a15-1964-4c43-b589-d1f0cc6bc596.e059ce618015e54d813a2480229315a514275c36.de-de.xlf"
$wsDe.Range("J3").Value = "01a97a15-1964-4c43-b589-d1f0cc6bc596.e059ce618015e54d813a2480229315a514275c36.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-30 07:06:35"
$wsDe.Range("K3").Value = "2016-08-30 07:06:35"

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($baseUrl + "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"), "", "", "01a97a15-1964-4c43-b589-d1f0cc6bc596.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($baseUrl + "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"), "", "", "01a97a15-1964-4c43-b589-d1f0cc6bc596.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($baseUrl + "ffffa6b67d8f-8e16-4a84-8d2f-d01ff423c2eb.md"), "", "", "ffffa6b67d8f-8e16-4a84-8d2f-d01ff423c2eb.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($baseUrl + "01a97a15-1964-4c43-b589-d1f0cc6bc596.md"), "", "", "01a97a15-1964-4c43-b589-d1f0cc6bc596.md")

$wsDe.Range("A2").Style = "HyperLink"
$wsDe.Range("A3").Style = "HyperLink"
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Range("I3").Style = "HyperLink"

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
